$wb = $excel.ActiveWorkbook
$wsNew = $wb.Worksheets.Item("new_datasets")
$wsExcl = $wb.Worksheets.Item("new_dataset_exclude")

# Rows 13:15 of "new_dataset_exclude" are being reclassified as "new_datasets"
# (figure re-ran and these three bioprojects moved out of the exclude list).
# Move them: copy formatting + values to the bottom of new_datasets, then
# delete the rows from new_dataset_exclude.
for ($i = 13; $i -le 15; $i++) {
    $srcRow = $wsExcl.Range("A" + $i + ":K" + $i)

    $srcRow.Copy()
    $wsNew.Range("A" + $i).PasteSpecial(-4122)

    for ($col = 1; $col -le 11; $col++) {
        $srcCell = $wsExcl.Cells.Item($i, $col)
        $dstCell = $wsNew.Cells.Item($i, $col)
        $dstCell.Value = $srcCell.Value2
    }

    $wsNew.Rows.Item($i).RowHeight = $wsExcl.Rows.Item($i).RowHeight
}

# Remove the now-duplicated rows from new_dataset_exclude
$wsExcl.Range("A13:K15").EntireRow.Delete()

# Restore selections to match what Excel would show after this edit
$wsExcl.Activate()
$wsExcl.Range("A13:K15").EntireRow.Select()

$wsNew.Activate()
$wsNew.Range("B15").Select()
